# Updates cryptos list values (Coin/Link/Price/Volume) per the scraped diff.
# Leading apostrophe forces text entry so values like "246.70" or "1.00" are not
# auto-coerced by Excel into numbers (which would drop formatting / trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.327.48"
$ws.Range("E2").Value = "'  +0.63%  "
$ws.Range("D3").Value = "'1.884.61"
$ws.Range("E3").Value = "'  -0.47%  "
$ws.Range("E4").Value = "'  -0.62%  "
$ws.Range("D5").Value = "'246.70"
$ws.Range("E5").Value = "'  -1.79%  "
$ws.Range("D6").Value = "'0.683"
$ws.Range("E6").Value = "'  -1.54%  "
$ws.Range("E7").Value = "'  -0.63%  "
$ws.Range("D8").Value = "'43.36"
$ws.Range("E8").Value = "'  +4.95%  "
$ws.Range("E9").Value = "'  +0.79%  "
$ws.Range("D10").Value = "'53.31"
$ws.Range("E10").Value = "'  +2.37%  "
$ws.Range("D11").Value = "'0.0749"
$ws.Range("E11").Value = "'  +0.31%  "
$ws.Range("D12").Value = "'0.0978"
$ws.Range("E12").Value = "'  +0.27%  "
$ws.Range("D13").Value = "'13.49"
$ws.Range("E13").Value = "'  +4.16%  "
$ws.Range("D14").Value = "'2.155.64"
$ws.Range("D15").Value = "'0.761"
$ws.Range("E15").Value = "'  +5.12%  "
$ws.Range("E16").Value = "'  -0.30%  "
$ws.Range("D17").Value = "'1.844.28"
$ws.Range("E17").Value = "'  -2.82%  "
$ws.Range("D18").Value = "'35.328.82"
$ws.Range("E18").Value = "'  +0.62%  "
$ws.Range("D19").Value = "'73.87"
$ws.Range("E19").Value = "'  +0.25%  "
$ws.Range("D20").Value = "'0.0₃0827"
$ws.Range("E20").Value = "'  -0.02%  "
$ws.Range("D21").Value = "'244.70"
$ws.Range("E21").Value = "'  -2.42%  "
$ws.Range("D22").Value = "'12.80"
$ws.Range("E22").Value = "'  -1.17%  "
$ws.Range("D23").Value = "'5.20"
$ws.Range("E23").Value = "'  +4.04%  "
$ws.Range("D24").Value = "'2.63"
$ws.Range("E24").Value = "'  +8.52%  "
$ws.Range("E25").Value = "'  -0.64%  "
$ws.Range("D26").Value = "'2.18"
$ws.Range("E26").Value = "'  -2.18%  "
$ws.Range("D27").Value = "'165.00"
$ws.Range("E27").Value = "'  -1.49%  "
$ws.Range("D28").Value = "'8.62"
$ws.Range("E28").Value = "'  +1.03%  "
$ws.Range("D29").Value = "'18.31"
$ws.Range("E29").Value = "'  -0.24%  "
$ws.Range("E30").Value = "'  -0.21%  "
$ws.Range("D31").Value = "'4.31"
$ws.Range("E31").Value = "'  +0.16%  "
$ws.Range("D32").Value = "'0.0596"
$ws.Range("E32").Value = "'  +1.13%  "
$ws.Range("E33").Value = "'  -0.47%  "
$ws.Range("B34").Value = "'WEMIXToken"
$ws.Range("C34").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'1.87"
$ws.Range("E34").Value = "'  -4.45%  "
$ws.Range("B35").Value = "'BinanceUSD"
$ws.Range("C35").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "'  -0.59%  "
$ws.Range("D36").Value = "'1.46"
$ws.Range("E36").Value = "'  -6.55%  "
$ws.Range("D37").Value = "'0.855"
$ws.Range("E37").Value = "'  +1.59%  "
$ws.Range("E38").Value = "'  -2.03%  "
$ws.Range("D39").Value = "'0.0728"
$ws.Range("E39").Value = "'  +9.79%  "
$ws.Range("D40").Value = "'17.61"
$ws.Range("E40").Value = "'  +0.02%  "
$ws.Range("E41").Value = "'  +2.64%  "
$ws.Range("D42").Value = "'97.23"
$ws.Range("E42").Value = "'  -1.17%  "
$ws.Range("E43").Value = "'  -1.56%  "
$ws.Range("E44").Value = "'  +0.99%  "
$ws.Range("D45").Value = "'1.310.65"
$ws.Range("E45").Value = "'  +1.10%  "
$ws.Range("D46").Value = "'0.0802"
$ws.Range("E46").Value = "'  +3.85%  "
$ws.Range("E47").Value = "'  -1.05%  "
$ws.Range("E48").Value = "'  -0.31%  "
$ws.Range("D49").Value = "'12.13"
$ws.Range("E49").Value = "'  +0.39%  "
$ws.Range("D50").Value = "'6.34"
$ws.Range("D51").Value = "'42.55"
$ws.Range("E51").Value = "'  +0.73%  "
